$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1955555555555556
$ws.Range("C2").Value = 0.5866666666666667
$ws.Range("J2").Value = 0.01777777777777778
$ws.Range("P2").Value = 0.1377777777777778
$ws.Range("S2").Value = 0.06222222222222222
$ws.Range("B3").Value = 0.01470588235294118
$ws.Range("C3").Value = 0.02205882352941177
$ws.Range("P3").Value = 0.7720588235294118
$ws.Range("S3").Value = 0.1911764705882353
$ws.Range("J4").Value = 0.07692307692307693
$ws.Range("P4").Value = 0.6153846153846154
$ws.Range("S4").Value = 0.3076923076923077
$ws.Range("B6").Value = 0.07222222222222222
$ws.Range("D6").Value = 0.02222222222222222
$ws.Range("F6").Value = 0.05
$ws.Range("J6").Value = 0.2388888888888889
$ws.Range("O6").Value = 0.01666666666666667
$ws.Range("Q6").Value = 0.1611111111111111
$ws.Range("R6").Value = 0.06111111111111111
$ws.Range("S6").Value = 0.3777777777777778
$ws.Range("B7").Value = 0.07792207792207792
$ws.Range("D7").Value = 0.01948051948051948
$ws.Range("F7").Value = 0.04545454545454546
$ws.Range("J7").Value = 0.1558441558441558
$ws.Range("Q7").Value = 0.2337662337662338
$ws.Range("R7").Value = 0.04545454545454546
$ws.Range("S7").Value = 0.4220779220779221
$ws.Range("B8").Value = 0.06471816283924843
$ws.Range("D8").Value = 0.01461377870563674
$ws.Range("F8").Value = 0.05845511482254697
$ws.Range("J8").Value = 0.104384133611691
$ws.Range("O8").Value = 0.00208768267223382
$ws.Range("Q8").Value = 0.208768267223382
$ws.Range("R8").Value = 0.07933194154488518
$ws.Range("S8").Value = 0.4676409185803758
$ws.Range("B9").Value = 0.06756756756756757
$ws.Range("D9").Value = 0.006756756756756757
$ws.Range("E9").Value = 0.006756756756756757
$ws.Range("F9").Value = 0.02702702702702703
$ws.Range("J9").Value = 0.1351351351351351
$ws.Range("Q9").Value = 0.1756756756756757
$ws.Range("R9").Value = 0.1013513513513514
$ws.Range("S9").Value = 0.4797297297297297
$ws.Range("B10").Value = 0.09430756159728122
$ws.Range("D10").Value = 0.02209005947323704
$ws.Range("E10").Value = 0.0008496176720475786
$ws.Range("F10").Value = 0.0713678844519966
$ws.Range("J10").Value = 0.1274426508071368
$ws.Range("O10").Value = 0.004248088360237893
$ws.Range("Q10").Value = 0.1954120645709431
$ws.Range("R10").Value = 0.1045029736618522
$ws.Range("S10").Value = 0.3797790994052676
$ws.Range("G11").Value = 0.1576763485477178
$ws.Range("J11").Value = 0.1203319502074689
$ws.Range("K11").Value = 0.2240663900414938
$ws.Range("L11").Value = 0.4854771784232365
$ws.Range("S11").Value = 0.01244813278008299
$ws.Range("G12").Value = 0.7916666666666666
$ws.Range("J12").Value = 0.125
$ws.Range("K12").Value = 0.008333333333333333
$ws.Range("L12").Value = 0.01666666666666667
$ws.Range("S12").Value = 0.05833333333333333
$ws.Range("G13").Value = 0.6097560975609756
$ws.Range("J13").Value = 0.2682926829268293
$ws.Range("S13").Value = 0.1219512195121951
$ws.Range("F15").Value = 0.01685393258426966
$ws.Range("H15").Value = 0.1404494382022472
$ws.Range("I15").Value = 0.07865168539325842
$ws.Range("J15").Value = 0.4325842696629214
$ws.Range("K15").Value = 0.07303370786516854
$ws.Range("M15").Value = 0.005617977528089887
$ws.Range("N15").Value = 0.005617977528089887
$ws.Range("O15").Value = 0.05617977528089887
$ws.Range("S15").Value = 0.1910112359550562
$ws.Range("F16").Value = 0.006578947368421052
$ws.Range("H16").Value = 0.1973684210526316
$ws.Range("I16").Value = 0.07894736842105263
$ws.Range("J16").Value = 0.4473684210526316
$ws.Range("K16").Value = 0.1052631578947368
$ws.Range("M16").Value = 0.05263157894736842
$ws.Range("O16").Value = 0.03947368421052631
$ws.Range("S16").Value = 0.07236842105263158
$ws.Range("F17").Value = 0.02392344497607655
$ws.Range("H17").Value = 0.2320574162679426
$ws.Range("I17").Value = 0.09090909090909091
$ws.Range("J17").Value = 0.3947368421052632
$ws.Range("K17").Value = 0.05980861244019139
$ws.Range("M17").Value = 0.0215311004784689
$ws.Range("O17").Value = 0.05980861244019139
$ws.Range("S17").Value = 0.1172248803827751
$ws.Range("F18").Value = 0.02072538860103627
$ws.Range("H18").Value = 0.1917098445595855
$ws.Range("I18").Value = 0.09844559585492228
$ws.Range("J18").Value = 0.4041450777202072
$ws.Range("K18").Value = 0.05699481865284974
$ws.Range("M18").Value = 0.03626943005181347
$ws.Range("O18").Value = 0.07253886010362694
$ws.Range("S18").Value = 0.1191709844559585
$ws.Range("F19").Value = 0.009136212624584718
$ws.Range("H19").Value = 0.2425249169435216
$ws.Range("I19").Value = 0.05398671096345515
$ws.Range("J19").Value = 0.3770764119601329
$ws.Range("K19").Value = 0.09883720930232558
$ws.Range("M19").Value = 0.0132890365448505
$ws.Range("O19").Value = 0.07142857142857142
$ws.Range("S19").Value = 0.1337209302325581
